$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow

# The sheet used to start with a blank placeholder row (old row 1) above the
# real header row (old row 2) and the single data row (old row 3). Drop that
# leading blank row so the header becomes row 1 and the data becomes row 2;
# this shifts every surviving cell (and its style) up by one row.
$ws.Rows.Item(1).Delete()

# Re-anchor the frozen panes now that the header moved from row 2 to row 1:
# freeze the header row together with the first data row (rows 1-2, column A)
# and have the scrollable area start at B3.
$win.FreezePanes = $false
$ws.Range("B3").Select()
$win.FreezePanes = $true

# Touch the very last row of the sheet so the worksheet's reported dimension
# covers the whole used column range through the last row of the grid, then
# clear the value back out again (only the row height sticks around).
$ws.Range("A1048576").Value = "x"
$ws.Range("A1048576").ClearContents()
$ws.Rows.Item(1048576).RowHeight = 12.8

# The header/first-data-row fill's (otherwise invisible) background tint
# changed from a grey to a blue tone; update it via the pattern colour while
# keeping the visible foreground colour untouched.
$ws.Range("A1:O2").Interior.PatternColor = 10040115
